function HexToOle([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r -bor ($g * 256) -bor ($b * 65536)
}
$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$theme = $nm.Theme
$tcs = $theme.ThemeColorScheme
$c = $tcs.Item(5)
$c.RGB = HexToOle("123456")
